$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3: Hydrogen / Non-metallic minerals value removed -> blank cell
$ws.Range("D3").Value = ""
$ws.Range("D3").Style = "Normal"

# C4: Methanol / Chemicals corrected to 0
$ws.Range("C4").Value = 0

# C5: Ammonia / Chemicals corrected to 0
$ws.Range("C5").Value = 0

# Row 7 relabelled from "Other" to "Biogas" with corrected value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 16.53329855345335

# New row 8 holds the (relabelled) "Other" category, pushed down by the
# Biogas row above it
$ws.Range("A8").Value = "Other"
$ws.Range("B8").Value = ""
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = ""
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = 29.90011565766796

# Match the row-label formatting (bold, centered, thin border) used by the
# other category rows in column A
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
